# Apply cryptocurrency price/volume updates per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    # Force the cell to store $Text as literal text even when it
    # looks like a number (e.g. "1.00"), then restore a plain/
    # unstyled "Normal" cell so no stray number-format style sticks.
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

# Row 2: D2: "68.207.02" -> "68.165.94"; E2: "  -3.61%  " -> "  -4.06%  "
$ws.Cells.Item(2, 4).Value = "68.165.94"
$ws.Cells.Item(2, 5).Value = "  -4.06%  "

# Row 3: D3: "3.690.29" -> "3.688.19"; E3: "  -4.38%  " -> "  -4.72%  "
$ws.Cells.Item(3, 4).Value = "3.688.19"
$ws.Cells.Item(3, 5).Value = "  -4.72%  "

# Row 4: D4: "0.999" -> "1.00"; E4: "  -0.18%  " -> "  +0.22%  "
Set-TextValue $ws.Cells.Item(4, 4) "1.00"
$ws.Cells.Item(4, 5).Value = "  +0.22%  "

# Row 5: D5: "593.84" -> "593.42"; E5: "  +0.35%  " -> "  -0.12%  "
Set-TextValue $ws.Cells.Item(5, 4) "593.42"
$ws.Cells.Item(5, 5).Value = "  -0.12%  "

# Row 6: D6: "181.70" -> "181.98"; E6: "  +8.97%  " -> "  +9.05%  "
Set-TextValue $ws.Cells.Item(6, 4) "181.98"
$ws.Cells.Item(6, 5).Value = "  +9.05%  "

# Row 7: D7: "3.677.86" -> "3.679.16"; E7: "  -4.60%  " -> "  -4.85%  "
$ws.Cells.Item(7, 4).Value = "3.679.16"
$ws.Cells.Item(7, 5).Value = "  -4.85%  "

# Row 8: E8: "  -6.58%  " -> "  -6.81%  "
$ws.Cells.Item(8, 5).Value = "  -6.81%  "

# Row 9: D9: "1.00" -> "0.997"; E9: "  -0.01%  " -> "  -0.24%  "
Set-TextValue $ws.Cells.Item(9, 4) "0.997"
$ws.Cells.Item(9, 5).Value = "  -0.24%  "

# Row 10: E10: "  -4.63%  " -> "  -5.22%  "
$ws.Cells.Item(10, 5).Value = "  -5.22%  "

# Row 11: E11: "  -7.11%  " -> "  -7.88%  "
$ws.Cells.Item(11, 5).Value = "  -7.88%  "

# Row 12: D12: "55.93" -> "55.92"; E12: "  +4.69%  " -> "  +4.55%  "
Set-TextValue $ws.Cells.Item(12, 4) "55.92"
$ws.Cells.Item(12, 5).Value = "  +4.55%  "

# Row 13: D13: "0.0000290" -> "0.0000289"; E13: "  -9.35%  " -> "  -9.98%  "
Set-TextValue $ws.Cells.Item(13, 4) "0.0000289"
$ws.Cells.Item(13, 5).Value = "  -9.98%  "

# Row 14: D14: "10.33" -> "10.32"; E14: "  -9.62%  " -> "  -9.88%  "
Set-TextValue $ws.Cells.Item(14, 4) "10.32"
$ws.Cells.Item(14, 5).Value = "  -9.88%  "

# Row 15: D15: "4.257.17" -> "4.272.88"; E15: "  -5.12%  " -> "  -4.86%  "
$ws.Cells.Item(15, 4).Value = "4.272.88"
$ws.Cells.Item(15, 5).Value = "  -4.86%  "

# Row 16: D16: "3.684.37" -> "3.686.85"; E16: "  -5.12%  " -> "  -4.96%  "
$ws.Cells.Item(16, 4).Value = "3.686.85"
$ws.Cells.Item(16, 5).Value = "  -4.96%  "

# Row 17: D17: "19.28" -> "19.25"; E17: "  -8.58%  " -> "  -9.76%  "
Set-TextValue $ws.Cells.Item(17, 4) "19.25"
$ws.Cells.Item(17, 5).Value = "  -9.76%  "

# Row 18: E18: "  -2.34%  " -> "  -2.42%  "
$ws.Cells.Item(18, 5).Value = "  -2.42%  "

# Row 19: E19: "  -7.02%  " -> "  -7.39%  "
$ws.Cells.Item(19, 5).Value = "  -7.39%  "

# Row 20: D20: "12.76" -> "12.75"; E20: "  -7.14%  " -> "  -7.87%  "
Set-TextValue $ws.Cells.Item(20, 4) "12.75"
$ws.Cells.Item(20, 5).Value = "  -7.87%  "

# Row 21: D21: "67.856.35" -> "67.970.09"; E21: "  -4.19%  " -> "  -4.05%  "
$ws.Cells.Item(21, 4).Value = "67.970.09"
$ws.Cells.Item(21, 5).Value = "  -4.05%  "

# Row 22: D22: "408.72" -> "408.08"; E22: "  -6.21%  " -> "  -6.66%  "
Set-TextValue $ws.Cells.Item(22, 4) "408.08"
$ws.Cells.Item(22, 5).Value = "  -6.66%  "

# Row 23: D23: "4.53" -> "4.54"; E23: "  -3.44%  " -> "  -3.74%  "
Set-TextValue $ws.Cells.Item(23, 4) "4.54"
$ws.Cells.Item(23, 5).Value = "  -3.74%  "

# Row 24: D24: "88.54" -> "88.55"; E24: "  -5.98%  " -> "  -6.27%  "
Set-TextValue $ws.Cells.Item(24, 4) "88.55"
$ws.Cells.Item(24, 5).Value = "  -6.27%  "

# Row 25: E25: "  -7.81%  " -> "  -8.13%  "
$ws.Cells.Item(25, 5).Value = "  -8.13%  "

# Row 26: E26: "  -7.88%  " -> "  -8.21%  "
$ws.Cells.Item(26, 5).Value = "  -8.21%  "

# Row 27: D27: "10.88" -> "10.89"; E27: "  -3.61%  " -> "  -4.25%  "
Set-TextValue $ws.Cells.Item(27, 4) "10.89"
$ws.Cells.Item(27, 5).Value = "  -4.25%  "

# Row 28: D28: "3.85" -> "3.84"; E28: "  -6.00%  " -> "  -5.94%  "
Set-TextValue $ws.Cells.Item(28, 4) "3.84"
$ws.Cells.Item(28, 5).Value = "  -5.94%  "

# Row 29: D29: "6.05" -> "6.03"; E29: "  +2.06%  " -> "  +1.83%  "
Set-TextValue $ws.Cells.Item(29, 4) "6.03"
$ws.Cells.Item(29, 5).Value = "  +1.83%  "

# Row 30: D30: "9.40" -> "9.39"; E30: "  -9.09%  " -> "  -9.72%  "
Set-TextValue $ws.Cells.Item(30, 4) "9.39"
$ws.Cells.Item(30, 5).Value = "  -9.72%  "

# Row 31: E31: "  -6.77%  " -> "  -7.13%  "
$ws.Cells.Item(31, 5).Value = "  -7.13%  "

# Row 32: D32: "7.27" -> "7.25"; E32: "  -11.40%  " -> "  -11.75%  "
Set-TextValue $ws.Cells.Item(32, 4) "7.25"
$ws.Cells.Item(32, 5).Value = "  -11.75%  "

# Row 33: D33: "12.43" -> "12.42"; E33: "  -8.26%  " -> "  -8.54%  "
Set-TextValue $ws.Cells.Item(33, 4) "12.42"
$ws.Cells.Item(33, 5).Value = "  -8.54%  "

# Row 34: E34: "  -7.00%  " -> "  -7.22%  "
$ws.Cells.Item(34, 5).Value = "  -7.22%  "

# Row 35: D35: "43.49" -> "43.39"; E35: "  -10.58%  " -> "  -10.10%  "
Set-TextValue $ws.Cells.Item(35, 4) "43.39"
$ws.Cells.Item(35, 5).Value = "  -10.10%  "

# Row 36: D36: "64.04" -> "64.15"; E36: "  -7.76%  " -> "  -8.51%  "
Set-TextValue $ws.Cells.Item(36, 4) "64.15"
$ws.Cells.Item(36, 5).Value = "  -8.51%  "

# Row 37: D37: "598.31" -> "599.36"; E37: "  -4.97%  " -> "  -5.06%  "
Set-TextValue $ws.Cells.Item(37, 4) "599.36"
$ws.Cells.Item(37, 5).Value = "  -5.06%  "

# Row 38: D38: "0.0₃0884" -> "0.0₃0881"; E38: "  -9.99%  " -> "  -10.49%  "
$ws.Cells.Item(38, 4).Value = "0.0₃0881"
$ws.Cells.Item(38, 5).Value = "  -10.49%  "

# Row 39: E39: "  +0.10%  " -> "  -0.01%  "
$ws.Cells.Item(39, 5).Value = "  -0.01%  "

# Row 40: E40: "  -6.86%  " -> "  -7.00%  "
$ws.Cells.Item(40, 5).Value = "  -7.00%  "

# Row 41: E41: "  -0.21%  " -> "  +0.17%  "
$ws.Cells.Item(41, 5).Value = "  +0.17%  "

# Row 42: E42: "  -7.87%  " -> "  -7.71%  "
$ws.Cells.Item(42, 5).Value = "  -7.71%  "

# Row 43: D43: "2.80" -> "2.78"; E43: "  +3.13%  " -> "  +2.79%  "
Set-TextValue $ws.Cells.Item(43, 4) "2.78"
$ws.Cells.Item(43, 5).Value = "  +2.79%  "

# Row 44: E44: "  -8.73%  " -> "  -9.30%  "
$ws.Cells.Item(44, 5).Value = "  -9.30%  "

# Row 45: E45: "  -7.13%  " -> "  -7.46%  "
$ws.Cells.Item(45, 5).Value = "  -7.46%  "

# Row 46: D46: "2.87" -> "2.86"; E46: "  -10.95%  " -> "  -12.55%  "
Set-TextValue $ws.Cells.Item(46, 4) "2.86"
$ws.Cells.Item(46, 5).Value = "  -12.55%  "

# Row 47: D47: "9.16" -> "9.17"; E47: "  -8.99%  " -> "  -9.20%  "
Set-TextValue $ws.Cells.Item(47, 4) "9.17"
$ws.Cells.Item(47, 5).Value = "  -9.20%  "

# Row 48: E48: "  -3.92%  " -> "  -3.65%  "
$ws.Cells.Item(48, 5).Value = "  -3.65%  "

# Row 49: B49: "ApeXProtocol" -> "Stellar"; C49: "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex" -> "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D49: "3.18" -> "0.134"; E49: "  -5.89%  " -> "  -7.16%  "
$ws.Cells.Item(49, 2).Value = "Stellar"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Cells.Item(49, 4) "0.134"
$ws.Cells.Item(49, 5).Value = "  -7.16%  "

# Row 50: B50: "Stellar" -> "ApeXProtocol"; C50: "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm" -> "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"; D50: "0.134" -> "3.16"; E50: "  -6.81%  " -> "  -5.69%  "
$ws.Cells.Item(50, 2).Value = "ApeXProtocol"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws.Cells.Item(50, 4) "3.16"
$ws.Cells.Item(50, 5).Value = "  -5.69%  "

# Row 51: D51: "2.733.24" -> "2.731.72"; E51: "  -3.65%  " -> "  -3.87%  "
$ws.Cells.Item(51, 4).Value = "2.731.72"
$ws.Cells.Item(51, 5).Value = "  -3.87%  "
